# Auto-generated edit script: updates cryptocurrency price/volume table
# to reflect refreshed values from the Wed Aug 21 21:32:50 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.869.47'
$ws.Range('E2').Value = '  +3.11%  '
$ws.Range('D3').Value = '2.610.73'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.26'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.54'
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  +1.13%  '
$ws.Range('D9').Value = '2.635.73'
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.53'
$ws.Range('E10').Value = '  -2.56%  '
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('E12').Value = '  -2.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.369'
$ws.Range('E13').Value = '  +7.14%  '
$ws.Range('D14').Value = '3.073.01'
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('D15').Value = '60.865.52'
$ws.Range('E15').Value = '  +3.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.48'
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('E17').Value = '  +3.10%  '
$ws.Range('D18').Value = '2.626.48'
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.30'
$ws.Range('E19').Value = '  +10.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.66'
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '347.86'
$ws.Range('E21').Value = '  +3.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.13'
$ws.Range('E22').Value = '  +14.17%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.523'
$ws.Range('E24').Value = '  +14.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.04'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.993'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.71'
$ws.Range('E28').Value = '  +6.39%  '
$ws.Range('D29').Value = '0.0₃0798'
$ws.Range('E29').Value = '  +2.64%  '
$ws.Range('E30').Value = '  +7.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.33'
$ws.Range('E32').Value = '  +4.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '161.23'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.52'
$ws.Range('E34').Value = '  +2.83%  '
$ws.Range('E35').Value = '  +5.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.966'
$ws.Range('E36').Value = '  +10.91%  '
$ws.Range('E37').Value = '  +4.78%  '
$ws.Range('E38').Value = '  +5.99%  '
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('E41').Value = '  +3.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '297.58'
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '139.50'
$ws.Range('E43').Value = '  +11.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0989'
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.995'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.605'
$ws.Range('E46').Value = '  +2.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0551'
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0241'
$ws.Range('E48').Value = '  +4.00%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.78'
$ws.Range('E49').Value = '  +7.19%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.69'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('D51').Value = '2.046.42'
$ws.Range('E51').Value = '  +5.11%  '
